# Apply weekly update: insert two new price rows for "Ají" at the top of the
# data block (rows 595-596), pushing the existing rows 595-638 down to 597-640.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 595.
$ws.Rows.Item(595).Insert()
$ws.Rows.Item(595).Insert()

# Copy the number style used by column D (date) from the row right below
# (now row 597, the original row 595) onto the two newly inserted D cells.
$ws.Cells.Item(597, 4).Copy()
$ws.Cells.Item(595, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(596, 4).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

function Set-AjiRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidades
    )

    $ws.Cells.Item($Row, 1).Value = 10
    $ws.Cells.Item($Row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($Row, 3).Value = "La Araucanía"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 9
    $ws.Cells.Item($Row, 6).Value = 100112021
    $ws.Cells.Item($Row, 7).Value = "Ají"
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = $Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# New row 595
Set-AjiRow 595 44610 "Americana (o)" "Primera" 160 13000 15000 14000 "`$/caja 15 kilos" "Región del Maule" 933 15

# New row 596
Set-AjiRow 596 44610 "Chilena(o)" "Primera" 40 15000 15000 15000 "`$/caja 15 kilos" "Región del Maule" 1000 15
